# Update "想去人数" (want-to-go count) figures in column F on both the
# "展览" and "全部类型" worksheets to reflect the newly generated output.

$wb = $excel.ActiveWorkbook

# Sheet "展览" (exhibitions)
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F4").Value = 1290
$wsExpo.Range("F16").Value = 613
$wsExpo.Range("F17").Value = 104
$wsExpo.Range("F24").Value = 2746
$wsExpo.Range("F37").Value = 315

# Sheet "全部类型" (all types)
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 1290
$wsAll.Range("F17").Value = 613
$wsAll.Range("F18").Value = 104
$wsAll.Range("F25").Value = 2746
$wsAll.Range("F38").Value = 315
